# Contest 23 CSK vs SRH
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the last contestant's contest/name (column AB/AC) from
# "Rag Nat0112" / "Raghu" to "Venni 3022" / "Venni"
$ws.Range("AB11").Value = "Venni 3022"
$ws.Range("AB12").Value = "Venni"

# Fill in the scores for match 23 (row 35, CSK vs SRH) for every player
$ws.Range("E35").Value = 0
$ws.Range("H35").Value = 60
$ws.Range("K35").Value = 30
$ws.Range("N35").Value = 50
$ws.Range("Q35").Value = 70
$ws.Range("T35").Value = 20
$ws.Range("W35").Value = 40
$ws.Range("Z35").Value = 100
$ws.Range("AC35").Value = 80
